# Generate Report for Handback
#
# The localization-status report regenerates with the two tracked files
# ("a5c21f88-...md" and "0de1a17f-...md") swapping their row order (the
# newest handback now sorts first), and "0de1a17f-...md" moves from
# "Ready for handoff" to "Handed back: in sync with en-US" with a fresh
# handback timestamp recorded for both locales.

$wb = $excel.ActiveWorkbook

$a5 = "a5c21f88-c3ef-43f4-a57a-934f0c15deac"
$de0 = "0de1a17f-000a-45e9-a2fd-87fe89ec20f2"

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$de0.md"
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack

$wsOverview.Range("A3").Value = "$a5.md"
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address(0, 0)
    if ($addr -eq "A2") {
        $hl.TextToDisplay = "$de0.md"
    } elseif ($addr -eq "A3") {
        $hl.TextToDisplay = "$a5.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 now carries the 0de1a17f file's data (handed back, new timestamp)
$wsZh.Range("A2").Value = "$de0.md"
$wsZh.Range("B2").Value = $handedBack
$wsZh.Range("C2").Value = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-07 09:54:56"
$wsZh.Range("E2").Value = "$de0.md"
$wsZh.Range("F2").Value = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-03-07 09:55:45"
$wsZh.Range("H2").Value = "Include"

# Row 3 now carries the a5c21f88 file's data (unchanged timestamps)
$wsZh.Range("A3").Value = "$a5.md"
$wsZh.Range("B3").Value = $handedBack
$wsZh.Range("C3").Value = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-07 09:52:41"
$wsZh.Range("E3").Value = "$a5.md"
$wsZh.Range("F3").Value = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-03-07 09:53:49"
$wsZh.Range("H3").Value = "Include"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address(0, 0)
    if ($addr -eq "A2") {
        $hl.TextToDisplay = "$de0.md"
    } elseif ($addr -eq "C2") {
        $hl.TextToDisplay = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.zh-cn.xlf"
    } elseif ($addr -eq "E2") {
        $hl.TextToDisplay = "$de0.md"
    } elseif ($addr -eq "F2") {
        $hl.TextToDisplay = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.zh-cn.xlf"
    } elseif ($addr -eq "A3") {
        $hl.TextToDisplay = "$a5.md"
    } elseif ($addr -eq "C3") {
        $hl.TextToDisplay = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.zh-cn.xlf"
    } elseif ($addr -eq "E3") {
        $hl.TextToDisplay = "$a5.md"
    } elseif ($addr -eq "F3") {
        $hl.TextToDisplay = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 now carries the 0de1a17f file's data (handed back, new timestamp)
$wsDe.Range("A2").Value = "$de0.md"
$wsDe.Range("B2").Value = $handedBack
$wsDe.Range("C2").Value = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-07 09:55:08"
$wsDe.Range("E2").Value = "$de0.md"
$wsDe.Range("F2").Value = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.de-de.xlf"
$wsDe.Range("G2").Value = "2016-03-07 09:56:07"
$wsDe.Range("H2").Value = "Include"

# Row 3 now carries the a5c21f88 file's data (unchanged timestamps)
$wsDe.Range("A3").Value = "$a5.md"
$wsDe.Range("B3").Value = $handedBack
$wsDe.Range("C3").Value = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-07 09:53:02"
$wsDe.Range("E3").Value = "$a5.md"
$wsDe.Range("F3").Value = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.de-de.xlf"
$wsDe.Range("G3").Value = "2016-03-07 09:54:11"
$wsDe.Range("H3").Value = "Include"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address(0, 0)
    if ($addr -eq "A2") {
        $hl.TextToDisplay = "$de0.md"
    } elseif ($addr -eq "C2") {
        $hl.TextToDisplay = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.de-de.xlf"
    } elseif ($addr -eq "E2") {
        $hl.TextToDisplay = "$de0.md"
    } elseif ($addr -eq "F2") {
        $hl.TextToDisplay = "$de0.b368d0d230f7d7e459371504b4d2f023c7e4ec69.de-de.xlf"
    } elseif ($addr -eq "A3") {
        $hl.TextToDisplay = "$a5.md"
    } elseif ($addr -eq "C3") {
        $hl.TextToDisplay = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.de-de.xlf"
    } elseif ($addr -eq "E3") {
        $hl.TextToDisplay = "$a5.md"
    } elseif ($addr -eq "F3") {
        $hl.TextToDisplay = "$a5.8367c4a9143281e4556ca787cdecb559fd2e26f5.de-de.xlf"
    }
}
